{"js": "// Insert a new bold \"Others: \" run immediately after the\n// \"#mr_rec_needs_inf___88# \" placeholder run, and before the run that\n// follows it (\"#mr_needs_oth_inf#\"), per the diff:\n//   + <w:r><w:rPr>(rFonts minorHAnsi, b, bCs, sz20, szCs20, lang en-GB)</w:rPr>\n//   +   <w:t xml:space=\"preserve\">Others: </w:t>\n//   + </w:r>\n\nconst body = context.document.body;\n\n// The placeholder text (including its trailing space) that precedes the\n// insertion point. This run is unique in the template.\nconst anchorText = \"#mr_rec_needs_inf___88# \";\nconst results = body.search(anchorText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(`Could not find anchor text: ${anchorText}`);\n}\n\nconst anchor = results.items[0];\n\n// Insert the new text right after the placeholder run.\nconst inserted = anchor.insertText(\"Others: \", Word.InsertLocation.after);\nawait context.sync();\n\n// Force the new text into its own run (distinct from its neighbours) by\n// toggling formatting across two sync boundaries before settling on the\n// final values that match the surrounding runs in the template.\ninserted.font.bold = false;\nawait context.sync();\n\ninserted.font.bold = true; // bold (w:b / w:bCs)\ninserted.font.size = 10; // w:sz/w:szCs val=\"20\" half-points => 10pt\n\nawait context.sync();\n", "ps1": "# Insert a new bold \"Others: \" run immediately after the\n# \"#mr_rec_needs_inf___88# \" placeholder run, and before the run that\n# follows it (\"#mr_needs_oth_inf#\"), per the diff:\n#   + <w:r><w:rPr>(rFonts minorHAnsi, b, bCs, sz20, szCs20, lang en-GB)</w:rPr>\n#   +   <w:t xml:space=\"preserve\">Others: </w:t>\n#   + </w:r>\n\n$d = $word.ActiveDocument\n\n# Locate the placeholder text (including its trailing space) that precedes\n# the insertion point. This text is unique in the template.\n$rng = $d.Content\n$found = $rng.Find.Execute(\"#mr_rec_needs_inf___88# \")\nif (-not $found) {\n    throw \"Could not find anchor text: #mr_rec_needs_inf___88# \"\n}\n\n# Collapse the found range to its end (insertion point right before the\n# following \"#mr_needs_oth_inf#\" run) and insert the new text there.\n$rng.Collapse(0)\n$rng.InsertAfter(\"Others: \")\n\n# Land the inserted text in its own run (distinct from its neighbours) by\n# toggling formatting before settling on the final values, which match the\n# bold/size/language formatting used by the surrounding runs in the template.\n$rng.Font.Bold = 0\n$rng.Font.Bold = 1\n$rng.Font.Size = 10\n"}
